$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Turn the plain-text e-mail address "nanxu@umd.edu" into a real
#    hyperlink pointing at the new address "im.nan.xu@gmail.com".
#    (This mirrors the diff: the single run containing
#    "...Dr. Nan Xu (nanxu@umd.edu) with..." is split into
#    "...Dr. Nan Xu (" + [hyperlink]im.nan.xu@gmail.com[/hyperlink] + ") with...")
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("nanxu@umd.edu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'nanxu@umd.edu' in the document"
}

$h = $d.Hyperlinks.Add($rng, "mailto:im.nan.xu@gmail.com", "", "", "im.nan.xu@gmail.com")
# Keep the same small font size (10pt / sz=20) that the rest of the
# paragraph uses -- Hyperlinks.Add resets the run to the default size.
$h.Range.Font.Size = 10
$h.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# 2) Add an extra blank paragraph right after the "How to Apply" /
#    "...in the subject line." paragraph (there was already one blank
#    paragraph there; the diff adds a second one before it).
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("in the subject line.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'in the subject line.' in the document"
}
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()
